# Applies the "data up to 10th" edit described by the diff:
#  - back-fills American Samoa (column E) values for rows 90-94
#  - fills in the state data (columns B:BE, minus E) for rows 184-188
#    which previously only had the date label in column A
#  - appends 5 new date rows (189-193) for 06-10 Aug 2020, label only
#
# The sheet dimension and the shared-strings table are maintained
# automatically by the engine as cells are written, so no explicit
# bookkeeping is needed for those.

function Set-RowRange {
    # Writes a contiguous block of row values starting at column
    # $startCol in one shot via a 1xN COM SAFEARRAY, which is how
    # Excel.Range.Value assignment works for multi-cell ranges.
    # $null entries leave the corresponding cell untouched/empty.
    param($ws, $row, $startCol, $values)
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $c1 = $ws.Cells.Item($row, $startCol)
    $c2 = $ws.Cells.Item($row, $startCol + $n - 1)
    $r = $ws.Range($c1, $c2)
    $r.Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Back-fill American Samoa (column E) for rows 90-94 ---------------
$ws.Cells.Item(90, 5).Value = 0.11762681762682
$ws.Cells.Item(91, 5).Value = 0.12646103896104
$ws.Cells.Item(92, 5).Value = 0.13155363155363
$ws.Cells.Item(93, 5).Value = 0.11111111111111
$ws.Cells.Item(94, 5).Value = 0.11825396825397

# --- Fill in state data for rows 184-188 (column A already set) ------
Set-RowRange $ws 184 2 @(0.05236894125691,0.027296820263775,0.023707883202681,$null,0.035251877167344,0.029524414628486,0.031689652957092,0.028404419992739,0.040242198681735,0.027127422245024,0.025647258870375,0.023404295448173,0.051012584614222,0.022187693039031,0.022104335057325,0.030570597412681,0.025844844580634,0.022145672993655,0.0219444344051,0.021626940012233,0.025950580008197,0.032314647521774,0.02493575273284,0.028810480048201,0.029448289572038,0.02674208436212,0.025943413737761,0.072089842336078,0.025454921578072,0.038681672315368,0.024686759002406,0.028495888883199,0.027390050852909,0.028931902665337,0.029333685864008,0.03832693941528,0.031691278068148,0.03552089782811,0.027051858377798,0.023920706178092,0.029094975030333,0.028959594587441,0.061677299532058,0.034541327023524,0.024866913221975,0.028964315749178,0.023905930422732,0.023019080346051,0.029369627775427,0.023130893472645,0.054284537967376,0.032329148349437,0.024987811164242,0.023330206607692,0.024701901418871,0.026808160556592)
Set-RowRange $ws 185 2 @(0.036108811242657,0.01924870504288,0.032627958971089,$null,0.033328141674001,0.042117049413847,0.028291681467377,0.022796971002466,0.032325367121722,0.021222858422007,0.017844892118754,0.01671571619194,0.055785477316229,0.020856311701279,0.023141886180335,0.025043047508718,0.022974036243869,0.021748152830385,0.024207523809012,0.024267006473062,0.029470781368128,0.034933952085989,0.029213702345996,0.034465594597767,0.028439041001512,0.025063936350384,0.022465067929249,0.069819804257756,0.02184072948115,0.032782936453231,0.021960030010419,0.025291584947085,0.023319425307389,0.02540085332594,0.025446986378121,0.035591063573775,0.026893030055301,0.03401947774191,0.024575892557192,0.021074732766132,0.025382180084329,0.024992476532166,0.065342833350522,0.027613929365035,0.020176938772928,0.024190826617114,0.019883156807217,0.02047482220067,0.026270553057887,0.023372690200667,0.064196083881365,0.031758675577317,0.026038195570553,0.024124474606093,0.02668518106044,0.028510268466991)
Set-RowRange $ws 186 2 @(0.06617922077982299,0.051435514606237,0.054789881853143,$null,0.052778129826449,0.053361079190811,0.050297996237894,0.044881125997,0.043961589402526,0.044483837753232,0.046056891536807,0.048432083987933,0.080424631595192,0.028280885519004,0.041636114087595,0.045848795270633,0.038857715629375,0.04195813286358,0.040028208714948,0.04178931108208,0.049013317617789,0.042023594188839,0.039753715307278,0.041993013383992,0.044799138350952,0.044280949507505,0.051436240615274,0.082742829013413,0.054434062255368,0.055377106854036,0.0461620828335,0.045565935670163,0.050846760072638,0.050310763068619,0.046832544151351,0.048888377160433,0.046671084572184,0.047019639162032,0.047478597485682,0.048963505680026,0.04561320446433,0.047270534871523,0.075926587244494,0.049864049698483,0.047495888603687,0.048287195478809,0.048641027540853,0.041777881971538,0.052504351937442,0.0436960433466,0.08370038473903101,0.0466021112452,0.040377403160971,0.047393313767626,0.046167392869518,0.045568076364838)
Set-RowRange $ws 187 2 @(0.061739479894878,0.05048680317387,0.051082669591503,$null,0.051170584580203,0.0499043201962,0.051529979194248,0.042302330509041,0.043345888552425,0.038730234790756,0.050278463056559,0.052337469506437,0.05424041180993,0.032372024379864,0.054605993687232,0.053496766889258,0.048195664216433,0.054869699754222,0.051052363744886,0.052060349549642,0.058446316263342,0.044823857450383,0.040469566831846,0.04822302852209,0.044491169314099,0.042456309770218,0.049278663432596,0.069897542774327,0.053564975857376,0.054277754060251,0.045719121309255,0.049352310025413,0.048429201703434,0.046239181265101,0.035503229712163,0.049399782900573,0.044117588926235,0.046031109656222,0.050374109307751,0.051527766072214,0.046059731556432,0.045104210623362,0.075373199365838,0.04533026982912,0.047980353542917,0.048454421690499,0.051034744653065,0.04358574310825,0.05552458577094,0.039616772782713,0.075657172511673,0.050246213599359,0.042971261139078,0.046660206993787,0.049983817076068,0.04604805183569)
Set-RowRange $ws 188 2 @(0.048411905948143,0.05273512361938,0.045304544090924,$null,0.039899622743965,0.047369273342078,0.051554762613423,0.043118538209497,0.048769779436759,0.045257751787951,0.052320370378101,0.05486187329576,0.054027635964131,0.034035182088718,0.050466896451127,0.053649840627189,0.045441691448378,0.051892143172501,0.047223206102648,0.049497641541111,0.057181755523456,0.045473213706517,0.044694861151463,0.048064704411606,0.045174890946347,0.041391897309527,0.04771561135948,0.070039234811703,0.05302809835105,0.053535638280576,0.051612126992161,0.051697799740053,0.048084063232789,0.046726283648042,0.042563417670102,0.049804262623915,0.044923055551527,0.049606905661489,0.052622079781145,0.053677055449764,0.049906018446746,0.051651726443984,0.078131756525274,0.051393577113784,0.053274934654203,0.052745786818004,0.055376785902723,0.048953231632517,0.060898464118539,0.052214241521135,0.08424430430877,0.05378486453325,0.047068741204581,0.052667469824907,0.054394438434353,0.050223078270777)

# --- Append new date rows 189-193 (label only, data comes later) -----
$ws.Cells.Item(189, 1).Value = "06 08 2020"
$ws.Cells.Item(190, 1).Value = "07 08 2020"
$ws.Cells.Item(191, 1).Value = "08 08 2020"
$ws.Cells.Item(192, 1).Value = "09 08 2020"
$ws.Cells.Item(193, 1).Value = "10 08 2020"
